# Lean Ux Canvas v2.0
# Alteracao: coloquei o titulo na parte de usuario e cliente que estava faltando.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Add the missing title textbox ("Clientes e usuários") on slide 1, above
#    the "Clientes e usuários" canvas section.
# ---------------------------------------------------------------------------
$s = $p.Slides.Item(1)

$left   = 707921 / 12700
$top    = 3491141 / 12700
$width  = 2855269 / 12700
$height = 461665 / 12700

$titleBox = $s.Shapes.AddTextbox(1, $left, $top, $width, $height)
$titleBox.Name = "CaixaDeTexto 28"

$titleBox.Fill.Visible = $false
$titleBox.TextFrame.WordWrap = $false
$titleBox.TextFrame.AutoSize = 1

$tr = $titleBox.TextFrame.TextRange
$tr.Text = "Clientes e usuários"
$tr.LanguageID = "pt-BR"
$tr.Font.Name = "MV Boli"
$tr.Font.NameComplexScript = "MV Boli"
$tr.Font.Size = 24

# ---------------------------------------------------------------------------
# 2) Refresh the "today" date footer field (11/09/2020 -> 14/09/2020) across
#    the slide master and every slide layout.
# ---------------------------------------------------------------------------
function Update-DateFooter($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.HasTextFrame) {
            $t = $sh.TextFrame.TextRange.Text
            if ($t -eq "11/09/2020") {
                $sh.TextFrame.TextRange.Text = "14/09/2020"
            }
        }
    }
}

$master = $p.SlideMaster
Update-DateFooter $master.Shapes

$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DateFooter $layouts.Item($li).Shapes
}
